$wb = $excel.ActiveWorkbook

# --- 1) Rename the "Requested quantity" header on the existing sheets -------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2) Add the new "PO Forecast" sheet at the end of the workbook ---------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "PO Forecast"

# Copy the header formatting (bold, border, centered) from the Weekly sheet
# header row so the new sheet reuses the same cell style.
$wsWeekly.Range("A1:B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

# Copy the date-column formatting (custom date number format) from the
# Weekly sheet down the full length of the new "ds" column.
$wsWeekly.Range("A2:A3").Copy()
$newSheet.Range("A2:A37").PasteSpecial(-4122)

$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# --- 3) Fill in the forecast data rows --------------------------------------
$forecastRows = @(
    @(44934.99999999999, 100, 42.60468097152076, 158.4779961404936),
    @(44941.99999999999, 98, 42.80789384183405, 157.1098338968943),
    @(44948.99999999999, 96, 38.64753981182123, 154.1885103939358),
    @(44955.99999999999, 94, 37.73407173647633, 153.0228975369994),
    @(44969.99999999999, 90, 26.52478693946361, 145.766340539953),
    @(44976.99999999999, 88, 26.23362395462122, 150.4849769196296),
    @(44983.99999999999, 86, 24.53143879393068, 146.154449655493),
    @(44990.99999999999, 84, 21.82255424141303, 143.337497498818),
    @(44997.99999999999, 82, 28.09984959329807, 142.0464453630807),
    @(45011.99999999999, 78, 14.68249916409069, 136.2608899453421),
    @(45025.99999999999, 74, 17.58430827177134, 133.7134768300588),
    @(45039.99999999999, 69, 11.55016243827772, 129.3992918696584),
    @(45046.99999999999, 67, 11.08110087230729, 127.8968809333313),
    @(45060.99999999999, 63, 2.583497650285437, 123.9686836684358),
    @(45067.99999999999, 61, 3.460309645124944, 119.3355694610427),
    @(45081.99999999999, 57, 0.3308281052358593, 116.4636081144845),
    @(45088.99999999999, 55, -3.72333993023717, 113.8515147960396),
    @(45095.99999999999, 53, -7.724849222316562, 115.9760234877345),
    @(45102.99999999999, 51, -10.03422330623467, 112.6497233729906),
    @(45109.99999999999, 49, -13.37763723821949, 109.6162704118612),
    @(45116.99999999999, 47, -10.08343583187461, 108.2449595973598),
    @(45123.99999999999, 45, -12.59207141039837, 99.95308941793263),
    @(45130.99999999999, 43, -16.21799218201941, 100.2114419316825),
    @(45137.99999999999, 41, -15.04690889891358, 100.5938892359681),
    @(45144.99999999999, 39, -23.0475801665318, 94.43933482846211),
    @(45151.99999999999, 36, -20.57674117417246, 97.12809072035228),
    @(45158.99999999999, 34, -24.57282818876186, 96.27685229721028),
    @(45165.99999999999, 32, -28.08913474710234, 98.1229220518801),
    @(45172.99999999999, 30, -31.26414264606728, 90.24965721582917),
    @(45179.99999999999, 28, -27.16659494677742, 89.6399891694012),
    @(45186.99999999999, 26, -31.77957174469023, 88.41484389982149),
    @(45193.99999999999, 24, -38.92038228546719, 81.54965655452706),
    @(45200.99999999999, 22, -30.50381623177162, 79.93548315570959),
    @(45207.99999999999, 20, -37.19018177201698, 78.65454221357284),
    @(45214.99999999999, 18, -40.79036612730513, 79.17952954098898),
    @(45221.99999999999, 16, -48.52972939364304, 76.6260200159046)
)

$r = 2
foreach ($row in $forecastRows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Restore the originally active sheet/selection.
$wsWeekly.Activate()
[void]$wsWeekly.Range("A1").Select()
